$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-92 down to 55-93.
$ws.Rows(54).Insert()

# Populate the newly inserted row 54 with the new data record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 44740
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 100114002
$ws.Range("G54").Value = "Camote"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 40
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = 20000
$ws.Range("N54").Value = "$/malla 20 kilos"
$ws.Range("O54").Value = "Perú"
$ws.Range("P54").Value = 1000
$ws.Range("Q54").Value = 20
$ws.Range("R54").Value = "Hortaliza"
